# Apply price-list updates to the "TODAY" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODAY")

# Update product description text (shared across rows 21 and 82)
$ws.Range("E21").Value = "Μπανάνες® Κολομβιας (Ζυγιζόμενο) / Kgr"
$ws.Range("E82").Value = "Μπανάνες® Κολομβιας (Ζυγιζόμενο) / Kgr"

# Update retail prices (column G) per diff
$priceUpdates = @{
    "G2"  = 14.95
    "G4"  = 15.45
    "G5"  = 2.2
    "G8"  = 13.9
    "G9"  = 13.9
    "G10" = 13.9
    "G13" = 16.4
    "G17" = 14.2
    "G18" = 10.9
    "G19" = 2.2
    "G21" = 1.48
    "G23" = 14.8
    "G25" = 3.95
    "G29" = 14.9
    "G31" = 13.55
    "G34" = 13.4
    "G37" = 14.95
    "G38" = 15.2
    "G40" = 9.9
    "G41" = 16.2
    "G43" = 9.9
    "G44" = 15.2
    "G58" = 17.9
    "G64" = 9.9
    "G65" = 15.98
    "G66" = 2.6
    "G75" = 2.95
    "G76" = 3.95
    "G77" = 2.2
    "G78" = 2.6
    "G82" = 1.48
    "G83" = 2.2
}

foreach ($cell in $priceUpdates.Keys) {
    $ws.Range($cell).Value = $priceUpdates[$cell]
}
